# Updates cryptocurrency price/volume figures in columns D and E
# to match the latest scrape, preserving the existing inline-string
# (text) cell type for every updated value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    # Leading apostrophe forces Excel to store the value as text
    # (matching the original t="inlineStr" cells) instead of
    # auto-converting number-looking strings like "1.00" into a
    # numeric value. Resetting the style afterwards avoids leaving
    # a stray quote-prefix cell style behind.
    $Cell.Value = "'" + $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "58.479.61"
Set-TextValue $ws.Range("E2") "  +0.50%  "
Set-TextValue $ws.Range("D3") "2.524.51"
Set-TextValue $ws.Range("E3") "  +1.81%  "
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.00%  "
Set-TextValue $ws.Range("D5") "521.65"
Set-TextValue $ws.Range("E5") "  +0.50%  "
Set-TextValue $ws.Range("D6") "132.93"
Set-TextValue $ws.Range("E6") "  -1.22%  "
Set-TextValue $ws.Range("E8") "  +1.20%  "
Set-TextValue $ws.Range("D9") "2.524.42"
Set-TextValue $ws.Range("E9") "  +1.17%  "
Set-TextValue $ws.Range("D10") "0.0978"
Set-TextValue $ws.Range("E10") "  -1.03%  "
Set-TextValue $ws.Range("E11") "  -1.52%  "
Set-TextValue $ws.Range("E12") "  -2.83%  "
Set-TextValue $ws.Range("E13") "  -2.49%  "
Set-TextValue $ws.Range("D14") "2.970.69"
Set-TextValue $ws.Range("E14") "  +1.71%  "
Set-TextValue $ws.Range("D15") "58.444.87"
Set-TextValue $ws.Range("E15") "  +0.37%  "
Set-TextValue $ws.Range("D16") "22.10"
Set-TextValue $ws.Range("E16") "  -0.27%  "
Set-TextValue $ws.Range("E17") "  -0.28%  "
Set-TextValue $ws.Range("D18") "2.525.31"
Set-TextValue $ws.Range("E18") "  +1.50%  "
Set-TextValue $ws.Range("D19") "10.66"
Set-TextValue $ws.Range("E19") "  -0.05%  "
Set-TextValue $ws.Range("D20") "322.01"
Set-TextValue $ws.Range("E20") "  +0.30%  "
Set-TextValue $ws.Range("E21") "  -0.64%  "
Set-TextValue $ws.Range("D22") "6.15"
Set-TextValue $ws.Range("E22") "  +6.93%  "
Set-TextValue $ws.Range("E23") "  +0.18%  "
Set-TextValue $ws.Range("D24") "64.67"
Set-TextValue $ws.Range("E24") "  +0.66%  "
Set-TextValue $ws.Range("D26") "1.00"
Set-TextValue $ws.Range("E26") "  +0.46%  "
Set-TextValue $ws.Range("E27") "  -1.11%  "
Set-TextValue $ws.Range("E28") "  +0.15%  "
Set-TextValue $ws.Range("D29") "0.0₃0751"
Set-TextValue $ws.Range("E29") "  +0.39%  "
Set-TextValue $ws.Range("D30") "168.43"
Set-TextValue $ws.Range("E30") "  -0.74%  "
Set-TextValue $ws.Range("E31") "  +1.52%  "
Set-TextValue $ws.Range("D32") "6.29"
Set-TextValue $ws.Range("E32") "  -0.50%  "
Set-TextValue $ws.Range("E33") "  -0.20%  "
Set-TextValue $ws.Range("E34") "  +0.02%  "
Set-TextValue $ws.Range("D35") "1.00"
Set-TextValue $ws.Range("E35") "  +0.25%  "
Set-TextValue $ws.Range("D36") "18.17"
Set-TextValue $ws.Range("E36") "  +0.29%  "
Set-TextValue $ws.Range("E37") "  -6.69%  "
Set-TextValue $ws.Range("D38") "3.91"
Set-TextValue $ws.Range("E38") "  -2.93%  "
Set-TextValue $ws.Range("E39") "  +0.98%  "
Set-TextValue $ws.Range("D40") "36.48"
Set-TextValue $ws.Range("E40") "  -0.39%  "
Set-TextValue $ws.Range("D41") "0.771"
Set-TextValue $ws.Range("E41") "  -3.67%  "
Set-TextValue $ws.Range("D42") "276.24"
Set-TextValue $ws.Range("E42") "  +0.00%  "
Set-TextValue $ws.Range("E43") "  -0.12%  "
Set-TextValue $ws.Range("D44") "129.82"
Set-TextValue $ws.Range("E44") "  +4.81%  "
Set-TextValue $ws.Range("E45") "  -3.61%  "
Set-TextValue $ws.Range("E46") "  +0.23%  "
Set-TextValue $ws.Range("E47") "  +0.80%  "
Set-TextValue $ws.Range("D48") "0.0499"
Set-TextValue $ws.Range("E48") "  +1.70%  "
Set-TextValue $ws.Range("E49") "  -0.35%  "
Set-TextValue $ws.Range("E50") "  +0.14%  "
Set-TextValue $ws.Range("D51") "16.90"
Set-TextValue $ws.Range("E51") "  -0.84%  "
